$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# 1. Remove the "Project_1-1" worksheet entirely.
$wb.Worksheets.Item("Project_1-1").Delete() | Out-Null

# 2. Update the Mentor register row with a new mentor's info.
$ws1 = $wb.Worksheets.Item("Mentor")
$ws1.Range("A8").Value = "육멘토"
$ws1.Range("B8").Value = "mentor666@gmail.com"

# 3. Update the Mentee register row with a new mentee's info,
#    turning the e-mail cell into a mailto hyperlink.
$ws2 = $wb.Worksheets.Item("Mentee")
$ws2.Range("A8").Value = "육멘티"
$ws2.Hyperlinks.Add($ws2.Range("B8"), "mailto:mentee66@gmail.com") | Out-Null
$ws2.Range("B8").Value = "mentee66@gmail.com"

# 4. Restore per-sheet selections (Mentor first so Mentee ends up active/selected last).
$ws1.Range("E8").Select() | Out-Null
$ws2.Range("D8").Select() | Out-Null
